$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '''28.336.81', '  -0.21%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '''1.864.19', '  -1.56%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '''1.022', '  +0.71%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '''316.18', '  -0.21%  ')
    ,@(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '''1.017', '  +0.32%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '''0.5094', '  -1.50%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '''0.3958', '  +0.79%  ')
    ,@(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '''0.08341', '  -1.03%  ')
    ,@(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '''1.106', '  -1.99%  ')
    ,@(11, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '''41.84', '  -0.42%  ')
    ,@(12, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '''6.222', '  -1.08%  ')
    ,@(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '''20.40', '  -1.58%  ')
    ,@(14, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '''1.832.58', '  -2.53%  ')
    ,@(15, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '''1.017', '  +0.23%  ')
    ,@(16, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '''7.186', '  -1.87%  ')
    ,@(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '''0.00001104', '  -0.66%  ')
    ,@(18, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '''90.47', '  -1.23%  ')
    ,@(19, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '''0.06728', '  -0.24%  ')
    ,@(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '''17.64', '  -1.48%  ')
    ,@(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '''1.017', '  +0.29%  ')
    ,@(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '''5.946', '  -2.14%  ')
    ,@(23, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '''28.359.51', '  -0.39%  ')
    ,@(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '''11.11', '  -0.87%  ')
    ,@(25, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '''2.283', '  +0.63%  ')
    ,@(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '''161.62', '  +0.67%  ')
    ,@(27, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '''2.042.12', '  -2.51%  ')
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '''20.62', '  -0.76%  ')
    ,@(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '''2.352', '  -5.41%  ')
    ,@(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '''126.95', '  +0.06%  ')
    ,@(31, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '''0.1044', '  -1.62%  ')
    ,@(32, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '''1.031', '  -1.18%  ')
    ,@(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '''5.766', '  -1.30%  ')
    ,@(34, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '''3.626', '  -0.31%  ')
    ,@(35, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '''0.02419', '  -2.18%  ')
    ,@(36, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '''0.06457', '  -2.41%  ')
    ,@(37, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '''0.2182', '  -1.70%  ')
    ,@(38, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '''8.810', '  -9.05%  ')
    ,@(39, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '''1.265', '  +1.24%  ')
    ,@(40, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '''1.173', '  -2.79%  ')
    ,@(41, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '''0.6371', '  -2.36%  ')
    ,@(42, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '''4.980', '  -0.68%  ')
    ,@(43, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '''11.21', '  -1.34%  ')
    ,@(44, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '''0.6002', '  -2.35%  ')
    ,@(45, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '''12.96', '  -2.11%  ')
    ,@(46, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '''3.701', '  -0.11%  ')
    ,@(47, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '''1.216', '  -5.51%  ')
    ,@(48, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '''1.984', '  -2.21%  ')
    ,@(49, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '''121.79', '  +0.24%  ')
    ,@(50, 'EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '''1.200', '  -3.43%  ')
    ,@(51, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '''0.06831', '  -1.62%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
}
